# Weekly update: two new price records for Cilantro (Vega Central Mapocho de
# Santiago) are inserted at the top of the data block starting at row 513.
# Inserting whole rows pushes the existing rows 513:609 down to 515:611,
# preserving all of their values/formatting, and grows the used range from
# A1:R609 to A1:R611 - matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 513-514; everything below shifts down by two rows.
$ws.Rows("513:514").Insert()

# --- New row 513 ---
$ws.Cells.Item(513, 1).Value = 9
$ws.Cells.Item(513, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(513, 3).Value = "Metropolitana"
$ws.Cells.Item(513, 4).Value = 44694
$ws.Cells.Item(513, 5).Value = 13
$ws.Cells.Item(513, 6).Value = 100112040
$ws.Cells.Item(513, 7).Value = "Cilantro"
$ws.Cells.Item(513, 8).Value = "Sin especificar"
$ws.Cells.Item(513, 9).Value = "Primera"
$ws.Cells.Item(513, 10).Value = 43
$ws.Cells.Item(513, 11).Value = 5000
$ws.Cells.Item(513, 12).Value = 5000
$ws.Cells.Item(513, 13).Value = 5000
$ws.Cells.Item(513, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(513, 15).Value = "Región Metropolitana"
$ws.Cells.Item(513, 16).Value = 139
$ws.Cells.Item(513, 17).Value = 36
$ws.Cells.Item(513, 18).Value = "Hortaliza"

# --- New row 514 ---
$ws.Cells.Item(514, 1).Value = 9
$ws.Cells.Item(514, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(514, 3).Value = "Metropolitana"
$ws.Cells.Item(514, 4).Value = 44694
$ws.Cells.Item(514, 5).Value = 13
$ws.Cells.Item(514, 6).Value = 100112040
$ws.Cells.Item(514, 7).Value = "Cilantro"
$ws.Cells.Item(514, 8).Value = "Sin especificar"
$ws.Cells.Item(514, 9).Value = "Primera"
$ws.Cells.Item(514, 10).Value = 160
$ws.Cells.Item(514, 11).Value = 9000
$ws.Cells.Item(514, 12).Value = 10000
$ws.Cells.Item(514, 13).Value = 9500
$ws.Cells.Item(514, 14).Value = "`$/docena de atados"
$ws.Cells.Item(514, 15).Value = "Región Metropolitana"
$ws.Cells.Item(514, 16).Value = 3167
$ws.Cells.Item(514, 17).Value = 3
$ws.Cells.Item(514, 18).Value = "Hortaliza"
